$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6122863292694092
$ws.Range("B1").Value = 2.655575513839722
$ws.Range("C1").Value = 3.222060441970825
$ws.Range("D1").Value = 3.761053562164307
$ws.Range("E1").Value = 0.9502931833267212
